$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 75. This shifts the existing rows 75-118 down to
# 76-119 (carrying their values/formatting with them), matching the diff
# where every row from 76 to 119 now holds what used to be in the row above
# it (75 -> 118), and the sheet dimension grows from A1:R118 to A1:R119.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new weekly record.
$ws.Cells.Item(75, 1).Value = 6
$ws.Cells.Item(75, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(75, 3).Value = "Metropolitana"
$ws.Cells.Item(75, 4).Value = 44438
$ws.Cells.Item(75, 5).Value = 13
$ws.Cells.Item(75, 6).Value = 100112022
$ws.Cells.Item(75, 7).Value = "Arveja Verde"
$ws.Cells.Item(75, 8).Value = "Perfection"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 170
$ws.Cells.Item(75, 11).Value = 34000
$ws.Cells.Item(75, 12).Value = 34000
$ws.Cells.Item(75, 13).Value = 34000
$ws.Cells.Item(75, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(75, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(75, 16).Value = 1360
$ws.Cells.Item(75, 17).Value = 25
$ws.Cells.Item(75, 18).Value = "Hortaliza"
